# Roadmap v2.xlsx — "FixBug after mesh 21/10/2018"
#
# Semantic changes on the "Roadmap" sheet:
#   - H2 "May danh (Decision making don gian)" moves into I2
#     (the old I2 "Test va sua loi" duplicate is overwritten/removed,
#      H2 becomes empty and drops out of the sheet entirely).
#   - F3 "Hieu ung cho linh, xe va cong trinh gom: ..." moves into J3
#     (F3 becomes empty and drops out of the sheet entirely).
#   - I5 date label is corrected from "Thu 7 - 27-10-2018" to
#     "Thu 3 - 30-10-2018".
#   - A brand-new merged cell K5:L5 is added with the date label
#     "Thu 3 - 6-11-2018", center-aligned, no fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

# --- Move H2 -> I2 (cut keeps the source's style on the destination,
#     then Clear drops the now-empty H2 cell from the sheet entirely) ---
$ws.Range("H2").Cut($ws.Range("I2"))
$ws.Range("H2").Clear()

# --- Move F3 -> J3 ---
$ws.Range("F3").Cut($ws.Range("J3"))
$ws.Range("F3").Clear()

# --- Add the new K5:L5 merged date cell ---
$ws.Range("K5:L5").HorizontalAlignment = -4108
$ws.Range("K5").Value = "Thứ 3 - 6-11-2018"
$ws.Range("K5:L5").Merge()

# --- Fix the I5 date label ---
$ws.Range("I5").Value = "Thứ 3 - 30-10-2018"

# --- Match the new selection left behind by the edit ---
$ws.Range("G6").Select()
